$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Input Elements") to hold the new
# "Skip Generation" field, shifting the existing Input Elements / Action /
# Expected Result columns one place to the right.
$ws.Columns.Item(5).Insert()

# Header for the new column.
$ws.Range("E1").Value = "Skip Generation"

# Populate the new column for the existing rows; row 6 is intentionally
# left blank.
$ws.Range("E2").Value = "yes"
$ws.Range("E3").Value = "yes"
$ws.Range("E4").Value = "yes"
$ws.Range("E5").Value = "yes"

# The freshly inserted column copied its per-row formatting from its left
# neighbour (column D); re-apply the same look the other data columns use
# (copied from column F, the old column E) so the new column matches them.
$ws.Range("F2:F5").Copy()
$ws.Range("E2:E5").PasteSpecial(-4122)

# Reflect the edit in the sheet's selection.
$ws.Range("E2:E5").Select()
